# 11th july run| ||
#
# Rows whose IVA (N), ValorTotal (O) and MontanteMB (P) columns were still
# the placeholder numeric 0 get back-filled with the looked-up amounts:
#   N -> "28€"     (IVA)
#   O -> "6,44€"   (ValorTotal)
#   P -> "34,44€"  (MontanteMB, i.e. N + O)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,5,7,8,9,10,11,12,13,14,15,16,17,18,21,22,26,30,33,35,36,42,43,47,49,51,52,54,56,57,58,59,60,62,63,68,70,72)

foreach ($r in $rows) {
    $nCell = $ws.Cells.Item($r, 14)

    # "28€" has no decimal separator, so a bare Value assignment would be
    # auto-parsed as the number 28 with a currency display format (same as
    # typing it straight into Excel). Force text entry instead, then drop
    # the format back to the sheet default so no style index sticks to the
    # cell (matches the other two columns, which stay text on their own).
    $nCell.NumberFormat = "@"
    $nCell.Value = "28€"
    $nCell.Style = "Normal"

    $ws.Cells.Item($r, 15).Value = "6,44€"
    $ws.Cells.Item($r, 16).Value = "34,44€"
}
